$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

# Replace the exact (whole) text of a paragraph (by 1-based index) with new
# text. Because Find/Replace on a range spanning multiple runs (and any
# w:proofErr siblings) collapses the match into a single run, this both
# updates the text and "flattens" any pre-existing run / proofErr splits.
function Set-ParagraphText($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex).Range
    $r = $d.Range($p.Start, $p.End)
    $f = $r.Find
    $f.ClearFormatting()
    $f.Text = $oldText
    $f.Replacement.Text = $newText
    $ok = $f.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)
    if (-not $ok) {
        Write-Host "Set-ParagraphText: NOT FOUND in paragraph $paraIndex : $oldText"
    }
    return $ok
}

# Find `searchText` inside paragraph `paraIndex` and replace the matched
# range's text with `newText` (keeps it inside whatever run(s) currently
# occupy that span; does not by itself introduce a run break).
function Replace-InParagraph($paraIndex, $searchText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex).Range
    $r = $d.Range($p.Start, $p.End)
    $f = $r.Find
    $f.ClearFormatting()
    $f.Text = $searchText
    $ok = $f.Execute()
    if (-not $ok) {
        Write-Host "Replace-InParagraph: NOT FOUND in paragraph $paraIndex : $searchText"
        return
    }
    $r.Text = $newText
}

# Force a run boundary right after the text `searchText` (matched starting
# from the beginning of paragraph `paraIndex`) by briefly adding and then
# removing a bookmark around the match. Word / the OOXML writer keeps the
# run split that resulted from the bookmark insertion even once the
# bookmark itself is deleted.
function Split-After($paraIndex, $searchText) {
    $p = $d.Paragraphs.Item($paraIndex).Range
    $r = $d.Range($p.Start, $p.End)
    $f = $r.Find
    $f.ClearFormatting()
    $f.Text = $searchText
    $ok = $f.Execute()
    if (-not $ok) {
        Write-Host "Split-After: NOT FOUND in paragraph $paraIndex : $searchText"
        return
    }
    $d.Bookmarks.Add("TempSplitMark", $r)
    $d.Bookmarks.Item("TempSplitMark").Delete()
}

# ---------------------------------------------------------------------------
# 1) "Assignment #" + "4"  ->  single run "Assignment #4"
# ---------------------------------------------------------------------------
Set-ParagraphText 2 "Assignment #4" "Assignment #4"

# ---------------------------------------------------------------------------
# 2) "August 3" + ", 2021"  ->  single run "August 3, 2021"
# ---------------------------------------------------------------------------
Set-ParagraphText 3 "August 3, 2021" "August 3, 2021"

# ---------------------------------------------------------------------------
# 3) Surviving / highest correlation bullet
#    "Out of the surviving patients, serum_sodium and platelets have the
#     highest correlation."
#    -> "Out of the surviving patients, serum_sodium and " | "serum_creatine " | "have the highest correlation."
# ---------------------------------------------------------------------------
Set-ParagraphText 9 `
    "Out of the surviving patients, serum_sodium and platelets have the highest correlation." `
    "Out of the surviving patients, serum_sodium and platelets have the highest correlation."
Replace-InParagraph 9 "platelets " "serum_creatine "
Split-After 9 "Out of the surviving patients, serum_sodium and "
Split-After 9 "Out of the surviving patients, serum_sodium and serum_creatine "

# ---------------------------------------------------------------------------
# 4) Surviving / lowest correlation bullet
#    "Out of the surviving patients, serum_creatine and serum_sodium have the
#     lowest correlation."
#    -> "Out of the surviving patients, " | "platelets " | "and serum_sodium have the lowest correlation."
# ---------------------------------------------------------------------------
Set-ParagraphText 10 `
    "Out of the surviving patients, serum_creatine and serum_sodium have the lowest correlation." `
    "Out of the surviving patients, serum_creatine and serum_sodium have the lowest correlation."
Replace-InParagraph 10 "serum_creatine and serum_sodium have the lowest correlation." "platelets and serum_sodium have the lowest correlation."
Split-After 10 "Out of the surviving patients, "
Split-After 10 "Out of the surviving patients, platelets "

# ---------------------------------------------------------------------------
# 5) Deceased / highest correlation bullet -> single merged run (no proofErr)
# ---------------------------------------------------------------------------
Set-ParagraphText 11 `
    "Out of deceased patients, serum_sodium and creatinine_phosphokinase have the highest correlation." `
    "Out of deceased patients, serum_sodium and creatinine_phosphokinase have the highest correlation."

# ---------------------------------------------------------------------------
# 6) Deceased / lowest correlation bullet
#    "Out of deceased patients, serum_creatine and serum_sodium have the
#     lowest correlation."
#    -> "Out of deceased patients, serum_creatine and " | [bookmark _Hlk78737708: "platelets" | " "] | "have the lowest correlation."
# ---------------------------------------------------------------------------
Set-ParagraphText 12 `
    "Out of deceased patients, serum_creatine and serum_sodium have the lowest correlation." `
    "Out of deceased patients, serum_creatine and serum_sodium have the lowest correlation."
Replace-InParagraph 12 "serum_sodium" "platelets"
Split-After 12 "Out of deceased patients, serum_creatine and platelets"
Split-After 12 "Out of deceased patients, serum_creatine and platelets "
$p12 = $d.Paragraphs.Item(12).Range
$bm = $d.Range($p12.Start, $p12.End)
$bmFind = $bm.Find
$bmFind.ClearFormatting()
$bmFind.Text = "platelets "
$bmFind.Execute() | Out-Null
$d.Bookmarks.Add("_Hlk78737708", $bm)

# ---------------------------------------------------------------------------
# 7) "For both deceased and surviving patients, ..." paragraph rewrite
# ---------------------------------------------------------------------------
$oldP13Tail = "serum_creatine and serum_sodium have the lowest correlation. But the variables for the highest correlation were different between surviving patients and deceased patients."
Set-ParagraphText 13 `
    "For both deceased and surviving patients, $oldP13Tail" `
    "For both deceased and surviving patients, $oldP13Tail"

$newTail = "different features have the lowest correlation. For instance, for surviving patients, serum_sodium and serum creatine have the highest correlation, where as for deceased patients, serum_sodium and creatinine_phosphokinase have the highest correlation. For lowest correlation, surviving patients and deceased patients, platelets had the lowest correlation, but with serum_sodium for surviving patients and with serum_creatine for deceased patients. This was not surprising to me as I would expect there to be a difference in values for deceased versus surviving patients."
Replace-InParagraph 13 $oldP13Tail $newTail

Split-After 13 "For both deceased and surviving patients, "
Split-After 13 "For both deceased and surviving patients, different features have the"
Split-After 13 "For both deceased and surviving patients, different features have the lowest correlation. "
Split-After 13 "For both deceased and surviving patients, different features have the lowest correlation. For instance, for surviving patients, serum_sodium and serum creatine have the highest correlation, where as for deceased patients, serum_sodium and creatinine_phosphokinase have the highest correlation. For lowest correlation, surviving patients and deceased patients, platelets had the lowest correlation, but with serum_sodium for surviving patients and with serum_creatine for deceased patients. This was not surprising to me as I would expect there to be a difference in values for deceased versus surviving patients"

Write-Host "Done."
